# Refresh cryptos price (D) and 1h volume-change (E) figures to match the
# latest coinranking.com snapshot pulled by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.909.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.551.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.97'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.772.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.556.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.519'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.899.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0469'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.411.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  +1.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.970'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.526'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('E42').Value = '  +3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.49%  '
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.686.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +4.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0958'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.11%  '
